# Update the "Neural_Net" row (row 4, columns B:L) on each of the three
# sheets (RMSE, MAE, R2) with the new metric values.

$wb = $excel.ActiveWorkbook

# --- RMSE sheet ---
$rmse = $wb.Worksheets.Item("RMSE")
$rmseValues = @(
    0.0072911809899999,
    0.0059464994,
    0.01001020223,
    0.01304606891,
    0.0124645997199999,
    0.0058072649799999,
    0.00951356301,
    0.01424377242,
    0.0277416998,
    0.13147457055,
    0.16399405654
)
$rmseArr = New-Object 'object[,]' 1, $rmseValues.Length
for ($i = 0; $i -lt $rmseValues.Length; $i++) {
    $rmseArr[0, $i] = $rmseValues[$i]
}
$rmse.Range("B4:L4").Value = $rmseArr

# --- MAE sheet ---
$mae = $wb.Worksheets.Item("MAE")
$maeValues = @(
    0.06314497619999999,
    0.0558334497,
    0.07011169239999999,
    0.07853884949999999,
    0.0716078533,
    0.0561192323999999,
    0.06687778029999999,
    0.0818240603,
    0.1114395889999999,
    0.1867503716,
    0.2265786485999999
)
$maeArr = New-Object 'object[,]' 1, $maeValues.Length
for ($i = 0; $i -lt $maeValues.Length; $i++) {
    $maeArr[0, $i] = $maeValues[$i]
}
$mae.Range("B4:L4").Value = $maeArr

# --- R2 sheet ---
$r2 = $wb.Worksheets.Item("R2")
$r2Values = @(
    0.6348107628557372,
    0.6541267060110949,
    0.7021641158261673,
    0.3473622471480092,
    0.5241735294732981,
    0.6264358415638843,
    0.6347964834682219,
    0.439415798211953,
    0.4287613672892593,
    0.3695403458867211,
    0.0103625134926402
)
$r2Arr = New-Object 'object[,]' 1, $r2Values.Length
for ($i = 0; $i -lt $r2Values.Length; $i++) {
    $r2Arr[0, $i] = $r2Values[$i]
}
$r2.Range("B4:L4").Value = $r2Arr
